# Regenerate orders with updated distance/sizes.
#
# The underlying face/fixation "Condition", "Filename_Left", "Filename_Right",
# "Distance" and "Size" columns encode a distance token (D51/D64/D80) and a
# size token (S30) inside otherwise-identical strings
# (e.g. "Face02_D51_S30", "Face02_D51_S30_l.png", "D51", "S30").
# This pass renumbers those tokens:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# (S20 and S25 are left untouched.) The substitution is applied to every
# string cell in the affected columns, across the whole used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row()
$numRows = $usedRange.Rows.Count()
$lastRow = $firstRow + $numRows - 1

# Columns that carry the distance/size-encoded strings:
#   B = Condition, D = Filename_Left, E = Filename_Right,
#   H = Distance, J = Size
$cols = @("B", "D", "E", "H", "J")

function Convert-Token([string]$text) {
    $result = $text
    $result = $result.Replace("D51", "D55")
    $result = $result.Replace("D64", "D69")
    $result = $result.Replace("D80", "D86")
    $result = $result.Replace("S30", "S31")
    return $result
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        $current = $cell.Value()
        if ($current -ne $null -and $current -is [string]) {
            $updated = Convert-Token $current
            if ($updated -ne $current) {
                $cell.Value = $updated
            }
        }
    }
}
